$wb = $excel.ActiveWorkbook

# Data-driven cell value updates: sheet name, cell address, new value
$updates = @(
    @("ALC", "H11", 6525.125),
    @("ALC", "I11", 6525.125),
    @("ALC", "K11", 6525.125),
    @("ALC", "M11", -6385.125),
    @("ALC", "H21", 18375.867),
    @("ALC", "I21", 14833.692),
    @("ALC", "J21", 41400),
    @("ALC", "K21", 14833.692),
    @("ALC", "L21", 41400),
    @("ALC", "M21", -14365.692),
    @("ALC", "N21", -42336),
    @("ALC", "H23", 18375.867),
    @("ALC", "I23", 14833.692),
    @("ALC", "J23", 41400),
    @("ALC", "K23", 14833.692),
    @("ALC", "L23", 41400),
    @("ALC", "M23", -14599.692),
    @("ALC", "N23", -41868),
    @("ALC", "H34", 1433.3334),
    @("ALC", "I34", 1433.3334),
    @("ALC", "J34", 0),
    @("ALC", "K34", 1433.3334),
    @("ALC", "L34", 0),
    @("ALC", "M34", -1230.3334),
    @("ALC", "H36", 1433.3334),
    @("ALC", "I36", 1433.3334),
    @("ALC", "J36", 0),
    @("ALC", "K36", 1433.3334),
    @("ALC", "L36", 0),
    @("ALC", "M36", -718.3334),
    @("ALC", "H132", 266060.34),
    @("ALC", "I132", 280774.8),
    @("ALC", "K132", 842324.3999999999),
    @("ALC", "M132", -839794.3999999999),
    @("ALC", "H137", 32260044),
    @("ALC", "I137", 1252.7307),
    @("ALC", "K137", 3758.1921),
    @("ALC", "M137", -1208.1921),
    @("ARM", "H32", 5491.7915),
    @("ARM", "I32", 6010.3613),
    @("ARM", "K32", 6010.3613),
    @("ARM", "M32", -5723.3613),
    @("ARM", "H92", 26513.334),
    @("ARM", "J92", 26513.334),
    @("ARM", "L92", 26513.334),
    @("BSM", "H105", 1886.6666),
    @("BSM", "I105", 1868.8889),
    @("BSM", "J105", 1993.3334),
    @("BSM", "K105", 1868.8889),
    @("BSM", "L105", 1993.3334),
    @("BSM", "M105", -121.8888999999999),
    @("BSM", "N105", -5487.3334),
    @("BSM", "H140", 34639),
    @("BSM", "J140", 34639),
    @("BSM", "L140", 34639),
    @("BSM", "N140", -44999),
    @("CRP", "H31", 1318.9318),
    @("CRP", "I31", 1021.5185),
    @("CRP", "J31", 1791.2941),
    @("CRP", "K31", 1021.5185),
    @("CRP", "L31", 1791.2941),
    @("CRP", "M31", -726.5185),
    @("CRP", "N31", -2381.2941),
    @("CRP", "H34", 1318.9318),
    @("CRP", "I34", 1021.5185),
    @("CRP", "J34", 1791.2941),
    @("CRP", "K34", 1021.5185),
    @("CRP", "L34", 1791.2941),
    @("CRP", "M34", -819.5185),
    @("CRP", "N34", -2195.2941),
    @("CRP", "H94", 861.52),
    @("CRP", "I94", 503.375),
    @("CRP", "J94", 1030.0588),
    @("CRP", "K94", 503.375),
    @("CRP", "L94", 1030.0588),
    @("CRP", "M94", -52.375),
    @("CRP", "N94", -1932.0588),
    @("CUL", "H4", 14529806),
    @("CUL", "J4", 495),
    @("CUL", "L4", 1485),
    @("CUL", "H9", 1263.3334),
    @("CUL", "I9", 0),
    @("CUL", "J9", 1263.3334),
    @("CUL", "K9", 0),
    @("CUL", "L9", 3790.0002),
    @("CUL", "N9", -4238.0002),
    @("CUL", "H131", 4217.2666),
    @("CUL", "I131", 6727.5),
    @("CUL", "J131", 3304.4546),
    @("CUL", "K131", 20182.5),
    @("CUL", "L131", 9913.363799999999),
    @("CUL", "M131", -15142.5),
    @("CUL", "N131", -19993.3638),
    @("CUL", "H140", 2387.0625),
    @("CUL", "I140", 1874.4166),
    @("CUL", "J140", 3925),
    @("CUL", "K140", 5623.2498),
    @("CUL", "L140", 11775),
    @("CUL", "M140", -443.2497999999996),
    @("CUL", "N140", -22135),
    @("GSM", "H92", 11639.8),
    @("GSM", "J92", 11639.8),
    @("GSM", "L92", 11639.8),
    @("GSM", "N92", -15383.8),
    @("LTW", "H22", 689.4),
    @("LTW", "I22", 850),
    @("LTW", "J22", 649.25),
    @("LTW", "K22", 850),
    @("LTW", "L22", 649.25),
    @("LTW", "M22", -555),
    @("LTW", "N22", -1239.25),
    @("LTW", "H27", 689.4),
    @("LTW", "I27", 850),
    @("LTW", "J27", 649.25),
    @("LTW", "K27", 850),
    @("LTW", "L27", 649.25),
    @("LTW", "M27", -743),
    @("LTW", "N27", -863.25),
    @("LTW", "H68", 2130),
    @("LTW", "I68", 2016.6666),
    @("LTW", "J68", 2300),
    @("LTW", "K68", 2016.6666),
    @("LTW", "L68", 2300),
    @("LTW", "M68", -1267.6666),
    @("LTW", "N68", -3798),
    @("LTW", "H71", 2130),
    @("LTW", "I71", 2016.6666),
    @("LTW", "J71", 2300),
    @("LTW", "K71", 10083.333),
    @("LTW", "L71", 11500),
    @("LTW", "M71", -6339.333000000001),
    @("LTW", "N71", -18988),
    @("WVR", "H62", 32373.75),
    @("WVR", "I62", 34665),
    @("WVR", "J62", 25500),
    @("WVR", "K62", 34665),
    @("WVR", "L62", 25500),
    @("WVR", "M62", -34041),
    @("WVR", "N62", -26748),
    @("WVR", "H65", 32373.75),
    @("WVR", "I65", 34665),
    @("WVR", "J65", 25500),
    @("WVR", "K65", 173325),
    @("WVR", "L65", 127500),
    @("WVR", "M65", -170205),
    @("WVR", "N65", -133740),
    @("WVR", "H140", 22371.5),
    @("WVR", "J140", 22371.5),
    @("WVR", "L140", 22371.5),
    @("WVR", "N140", -32731.5),
    @("WVR", "H141", 62701.875),
    @("WVR", "J141", 62701.875),
    @("WVR", "L141", 62701.875),
    @("WVR", "N141", -73061.875),
    @("ARM", "N92", -31505.334),
    @("CUL", "N4", -1709)
)

foreach ($item in $updates) {
    $sheetName = $item[0]
    $addr = $item[1]
    $val = $item[2]
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range($addr).Value = $val
}

# Cells that must be cleared entirely (no longer present after the edit)
$clears = @(
    @("ALC", "N34"),
    @("ALC", "N36"),
    @("CUL", "M9")
)

foreach ($item in $clears) {
    $sheetName = $item[0]
    $addr = $item[1]
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range($addr).ClearContents()
}

Write-Host "Applied $($updates.Count) updates and $($clears.Count) clears."
